$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting all existing data right by one column.
$ws.Columns("A:A").Insert()

# Populate the new RefID column with the Jira reference IDs for each test case row.
$ws.Range("A1").Value = "RefID"
$ws.Range("A2").Value = "LATFLD-24"
$ws.Range("A3").Value = "LATFLD-23"
$ws.Range("A4").Value = "LATFLD-17"
$ws.Range("A5").Value = "LATFLD-12"

# Match the bold/shaded header formatting already used by the other header cells.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats

# Resize the columns whose content/width changed as a result of the new data.
$ws.Columns("A:A").ColumnWidth = 11.166666666666666
$ws.Columns("E:E").ColumnWidth = 31.166666666666668
$ws.Columns("G:G").ColumnWidth = 14.877604166666666

# Re-create the (hidden) _FilterDatabase name that Excel maintains for the sheet's filtered range.
$name = $ws.Names.Add("_xlnm._FilterDatabase", "=Transmittals_Close_Cancel!`$A`$1:`$Q`$5")
$name.Visible = $False
